$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7272.422255400007
$ws.Range("D2").Value = 353.5872912166667
$ws.Range("B3").Value = 6760.700009166674
$ws.Range("D3").Value = 328.8478758333333
$ws.Range("B4").Value = 7236.999968266674
$ws.Range("D4").Value = 349.7912609333333
$ws.Range("B5").Value = 7006.445791766673
$ws.Range("D5").Value = 331.3696418
$ws.Range("B6").Value = 7266.316580000007
$ws.Range("D6").Value = 353.7159776
$ws.Range("B7").Value = 6990.85414203334
$ws.Range("D7").Value = 334.4684250833333
$ws.Range("B8").Value = 7238.452864566674
$ws.Range("D8").Value = 347.1777969333334
$ws.Range("B9").Value = 7266.884610816674
$ws.Range("D9").Value = 355.2152107833333
$ws.Range("B10").Value = 6977.038504816675
$ws.Range("D10").Value = 344.4895291833333
$ws.Range("B11").Value = 7246.483091233341
$ws.Range("D11").Value = 350.16805895
$ws.Range("B12").Value = 6999.71887263334
$ws.Range("D12").Value = 341.6009272833333
$ws.Range("B13").Value = 6964.543750066674
$ws.Range("D13").Value = 329.9399408166667
